# Update the slide format:
#  1. Bump the fixed "date last edited" placeholder text that appears on the
#     notes master, the slide master and every slide layout from
#     06/02/2018 -> 07/02/2018.
#  2. Fix the capitalisation of the title on slide 1:
#     "Worked Example of a public function" -> "Worked example of a public function".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "06/02/2018") {
                $tr.Text = "07/02/2018"
            }
        }
    }
}

# Notes master.
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide 1 title text correction.
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Worked Example of a public function") {
            $tr.Text = "Worked example of a public function"
        }
    }
}
